# Update countries & provincias Spain
# Applies the diff: country-row resorts (new case counts causing a handful
# of rows to swap position) plus straightforward numeric refreshes, and the
# "last updated" timestamp string on A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: worksheet row number, optional new country name for column A
# (only set when the row's occupant changed because of a re-sort), and the
# new B..H values (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes).
$updates = @(
    @{ Row=7;   Name=$null;             B=1122241; C=6431; D=923699; E=178743; F=0; G=150; H=19799 }  # Rusia
    @{ Row=26;  Name=$null;             B=257388;  C=4465; D=187958; E=59453;  F=0; G=140; H=9977  }  # Indonesia
    @{ Row=27;  Name=$null;             B=200041;  C=6667; D=144401; E=54324;  F=0; G=31;  H=1316  }  # Israel
    @{ Row=67;  Name="Austria";         B=39984;   C=681;  D=30949;  E=8258;   F=0; G=6;   H=777   }  # Austria moves above Azerbaiyan
    @{ Row=68;  Name="Azerbaiyan";      B=39378;   C=0;    D=36949;  E=1851;   F=0; G=0;   H=578   }  # Azerbaiyan shifts down
    @{ Row=89;  Name=$null;             B=15340;   C=204;  D=13815;  E=1268;   F=0; G=2;   H=257   }  # Croacia
    @{ Row=110; Name="Eslovaquia";      B=7269;    C=338;  D=3888;   E=3340;   F=0; G=1;   H=41    }  # Eslovaquia moves above Mozambique/Birmania
    @{ Row=111; Name="Mozambique";      B=7114;    C=0;    D=4064;   E=3005;   F=0; G=0;   H=45    }  # Mozambique shifts down
    @{ Row=112; Name="Birmania";        B=6959;    C=216;  D=1951;   E=4892;   F=0; G=1;   H=116   }  # Birmania shifts down
    @{ Row=128; Name=$null;             B=4694;    C=136;  D=3168;   E=1383;   F=0; G=1;   H=143   }  # Eslovenia
    @{ Row=131; Name=$null;             B=4140;    C=227;  D=1643;   E=2472;   F=0; G=2;   H=25    }  # Georgia
    @{ Row=133; Name="Lituania";        B=3932;    C=73;   D=2246;   E=1599;   F=0; G=0;   H=87    }  # Lituania moves above Siria
    @{ Row=134; Name="Siria";           B=3877;    C=0;    D=983;    E=2716;   F=0; G=0;   H=178   }  # Siria shifts down
    @{ Row=143; Name="Estonia";         B=3033;    C=57;   D=2387;   E=582;    F=0; G=0;   H=64    }  # Estonia moves above Mali
    @{ Row=144; Name="Mali";            B=3030;    C=0;    D=2380;   E=521;    F=0; G=0;   H=129   }  # Mali shifts down
    @{ Row=214; Name="Islas Malvinas";  B=13;      C=0;    D=13;     E=0;      F=0; G=0;   H=0     }  # Islas Malvinas moves above Montserrat
    @{ Row=215; Name="Montserrat";      B=13;      C=0;    D=12;     E=0;      F=0; G=0;   H=1     }  # Montserrat shifts down
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.Name) {
        $ws.Cells.Item($r, 1).Value = $u.Name
    }
    $ws.Cells.Item($r, 2).Value = $u.B
    $ws.Cells.Item($r, 3).Value = $u.C
    $ws.Cells.Item($r, 4).Value = $u.D
    $ws.Cells.Item($r, 5).Value = $u.E
    $ws.Cells.Item($r, 6).Value = $u.F
    $ws.Cells.Item($r, 7).Value = $u.G
    $ws.Cells.Item($r, 8).Value = $u.H
}

# Refresh the "last updated" footer string (A1).
$ws.Range("A1").Value = "Datos actualizados a 23 de Septiembre de 2020 a las 10:51"
